{"js": "// \"Finished event table notes\" \u2014 remove the leftover full-page screenshot\n// (\"Picture 4\") that was floating (anchored, behindDoc) in the paragraph\n// right after the \"Research and Assumptions\" heading. The picture's host\n// paragraph is left in place (now empty), matching the target edit.\n//\n// The picture is a floating/anchored drawing (wp:anchor), so it shows up in\n// context.document.body.shapes (not body.inlinePictures, which only covers\n// inline wp:inline drawings such as the cover-page logo).\n\nconst shapes = context.document.body.shapes;\nshapes.load(\"items/name\");\nawait context.sync();\n\nfor (let i = 0; i < shapes.items.length; i++) {\n  const shape = shapes.items[i];\n  if (shape.name === \"Picture 4\") {\n    shape.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"Finished event table notes\" \u2014 remove the leftover full-page screenshot\n# (\"Picture 4\") that was floating (anchored, behindDoc) in the paragraph\n# right after the \"Research and Assumptions\" heading. The picture's host\n# paragraph is left in place (now empty), matching the target edit.\n#\n# The picture is a floating/anchored drawing, so it lives in\n# $d.Shapes (not $d.InlineShapes, which only covers inline pictures such as\n# the cover-page logo).\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Shapes.Count; $i -ge 1; $i--) {\n  $shape = $d.Shapes.Item($i)\n  if ($shape.Name -eq \"Picture 4\") {\n    $shape.Delete()\n  }\n}\n"}
